$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (G1) onto the new H1 header
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column data values for rows 2-5
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
